$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Heures TD" column (F) to hold
# the new "Groupes CM" data; everything from old F onward shifts right by one.
$ws.Columns("F:F").Insert()

# New header for the inserted column.
$ws.Range("F1").Value = "Groupes CM"

# New "Groupes CM" values for each data row (rows with double CM sessions get 3).
$cmGroups = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 3
    14 = 1
    15 = 1
    16 = 1
    17 = 1
}
foreach ($r in $cmGroups.Keys) {
    $ws.Cells.Item($r, 6).Value = $cmGroups[$r]
}

# Bug fix: rows for ALG044C1 / ALG044C2 (now rows 14 & 15) shared the same
# "Groupes TD" value (1) because of duplicate event IDs; correct it to 2
# (now column H after the column insert).
$ws.Cells.Item(14, 8).Value = 2
$ws.Cells.Item(15, 8).Value = 2

# Re-register the sort definition over the expanded range so the sheet's
# stored sort-state reference grows from A2:I35 to A2:J35 (no data movement,
# the sheet is already in sorted order).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B35"))
$ws.Sort.SortFields.Add($ws.Range("A2:A35"))
$ws.Sort.SetRange($ws.Range("A2:J35"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Move the active selection to A17, matching the edited workbook's saved view.
[void]$ws.Range("A17").Select()
